# Complain.xlsx update — add "Trạng thái" (Status) column before the
# existing "Bộ phận đang XL" column (old column H), per commit "update
# code from duongdx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at H — this shifts every column from H onward
#    one place to the right (H->I, I->J, ... O->P) and extends all the
#    row "spans", merged cells, the used dimension, etc. automatically.
$ws.Columns("H:H").Insert()

# 2) Give the new column a header in the header row (row 16), using the
#    new shared string "Trạng thái".
$ws.Range("H16").Value = "Trạng thái"

# 3) Re-establish the AutoFilter over the now-wider header row
#    (was A17:O17, needs to become A17:P17). Turning filtering off first
#    forces Excel to rebuild the filter range instead of just keeping the
#    stale one after the column insert.
$ws.AutoFilterMode = $false
$ws.Range("A17:P17").AutoFilter()

# 4) The workbook-level hidden _FilterDatabase defined name also needs to
#    track the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$17:`$P`$17"
    }
}

# 5) Resize the columns around the inserted one to the widths the author
#    ended up with (manual column resizing after the insert).
$ws.Columns("D:D").ColumnWidth = 17.333333333333332
$ws.Columns("E:E").ColumnWidth = 11.666666666666666
$ws.Columns("F:F").ColumnWidth = 19.666666666666668
$ws.Columns("H:H").ColumnWidth = 11.166666666666666
$ws.Columns("I:I").ColumnWidth = 13.666666666666666
$ws.Columns("J:J").ColumnWidth = 14.833333333333334
$ws.Columns("K:K").ColumnWidth = 10.666666666666666
$ws.Columns("L:L").ColumnWidth = 13.333333333333334
$ws.Columns("M:M").ColumnWidth = 11.333333333333334
$ws.Columns("N:N").ColumnWidth = 10.5
$ws.Columns("O:O").ColumnWidth = 12.833333333333334
$ws.Columns("P:P").ColumnWidth = 12.5

# 6) Row 17 is a sample/placeholder numbering row (1..15 across A17:O17).
#    The column insert shifted it like every other row (H17 blank,
#    I17..P17 = 8..15), but the author actually just continued the
#    sequence in place, so put 8..16 back across H17:P17.
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 11
$ws.Range("L17").Value = 12
$ws.Range("M17").Value = 13
$ws.Range("N17").Value = 14
$ws.Range("O17").Value = 15
$ws.Range("P17").Value = 16

# 7) Put the active selection on the new column's header cell, matching
#    where the author was working.
$ws.Range("I16").Select()
